$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a "plain number" string (e.g. "585.29")
# must be forced to Text so Excel does not auto-convert them to a numeric
# type (the source workbook stores every Price/Volume cell as text).
$numericPriceCells = @("D5", "D6", "D14", "D19", "D20", "D23", "D30", "D31", "D32", "D33", "D36", "D38", "D39", "D42", "D43", "D45", "D48")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.236.18"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "2.481.63"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "585.29"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").Value = "173.89"
$ws.Range("E6").Value = "  +3.57%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "2.925.77"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").Value = "25.49"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "67.136.78"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").Value = "2.476.45"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "10.97"
$ws.Range("E19").Value = "  -1.34%  "

$ws.Range("D20").Value = "350.80"
$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "69.03"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").Value = "2.608.11"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").Value = "0.0₃0911"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").Value = "505.83"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").Value = "7.76"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  +1.15%  "

$ws.Range("D36").Value = "161.18"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").Value = "18.19"
$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("D39").Value = "1.34"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  +1.91%  "

$ws.Range("D42").Value = "0.329"
$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").Value = "4.85"
$ws.Range("E43").Value = "  +1.19%  "

$ws.Range("E44").Value = "  +1.84%  "

$ws.Range("D45").Value = "142.81"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").Value = "0.0₆0261"
$ws.Range("E46").Value = "  +3.26%  "

$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "0.515"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("E51").Value = "  +0.70%  "

# Strip the temporary "@" number format back off again so the only
# persisted difference versus the original cells is their text content
# (matches the source, which applies no explicit style to these cells).
$numericPriceRange = $ws.Range($numericPriceCells[0])
foreach ($addr in $numericPriceCells) {
    $numericPriceRange = $excel.Union($numericPriceRange, $ws.Range($addr))
}
$numericPriceRange.ClearFormats()

